# Update column C ("Förändrad") date values from 2023-09-15 (45184) to
# 2023-09-16 (45185) for every data row (rows 2-15) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45185
    }
}
